$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision regarding which movie to acquire was not reached.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision,`" indicating that no consensus was achieved regarding which movie to show on Friday.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday resulted in no conclusion, and therefore, the conversation about the movie acquisition has ended without a choice.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to show on Friday.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision-making process concluded without a definitive choice for Friday's movie, so I have noted that no decision was made.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"

$ws.Range("D11").Value = "Barbie_was_selected, "
